# Refresh the cryptos table: latest Price (D) / Volume(1h) change (E) per coin,
# plus rows 48-49 swap (Monero <-> Arweave) with their refreshed figures.
# Numeric-looking text values are entered with a leading apostrophe (and the
# cell style reset to Normal afterwards) so they stay plain text, matching the
# source data which stores prices like "0.999" / "69.047.97" as strings, not numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.047.97"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +1.88%  "
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "3.749.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.03%  "
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'602.18"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.72%  "
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'168.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.00%  "
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "3.748.66"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.03%  "
$ws.Range("E7").Style = "Normal"
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("E8").Style = "Normal"
$ws.Range("E9").Value = "  +0.59%  "
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.166"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'6.43"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +3.28%  "
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.460"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.09%  "
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'38.06"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.07%  "
$ws.Range("E13").Style = "Normal"
$ws.Range("D14").Value = "'0.0000247"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +1.69%  "
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "4.369.14"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.81%  "
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "3.732.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.25%  "
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "69.079.49"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +1.92%  "
$ws.Range("E17").Style = "Normal"
$ws.Range("D18").Value = "'7.30"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "
$ws.Range("E18").Style = "Normal"
$ws.Range("D19").Value = "'0.113"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").Value = "'17.18"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.75%  "
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'10.85"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +17.60%  "
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'493.39"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +1.08%  "
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.725"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.05%  "
$ws.Range("E23").Style = "Normal"
$ws.Range("D24").Value = "'0.0000151"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +6.32%  "
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'84.76"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.10%  "
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'2.31"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.21%  "
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'12.31"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -0.13%  "
$ws.Range("E27").Style = "Normal"
$ws.Range("D28").Value = "'10.12"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.13%  "
$ws.Range("E28").Style = "Normal"
$ws.Range("E29").Value = "  -0.05%  "
$ws.Range("E29").Style = "Normal"
$ws.Range("D30").Value = "'2.98"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +1.99%  "
$ws.Range("E30").Style = "Normal"
$ws.Range("D31").Value = "'2.50"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.71%  "
$ws.Range("E31").Style = "Normal"
$ws.Range("D32").Value = "'8.06"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.45%  "
$ws.Range("E32").Style = "Normal"
$ws.Range("D33").Value = "'31.61"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.70%  "
$ws.Range("E33").Style = "Normal"
$ws.Range("D34").Value = "3.891.73"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.97%  "
$ws.Range("E34").Style = "Normal"
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E35").Style = "Normal"
$ws.Range("D36").Value = "3.682.29"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.67%  "
$ws.Range("E36").Style = "Normal"
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.17%  "
$ws.Range("E37").Style = "Normal"
$ws.Range("E38").Value = "  +1.18%  "
$ws.Range("E38").Style = "Normal"
$ws.Range("D39").Value = "'5.86"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("E39").Style = "Normal"
$ws.Range("E40").Value = "  +1.08%  "
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.324"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +0.24%  "
$ws.Range("E41").Style = "Normal"
$ws.Range("E42").Value = "  +3.95%  "
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'430.44"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.35%  "
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'48.57"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.52%  "
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'1.98"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.34%  "
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'8.49"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +0.59%  "
$ws.Range("E46").Style = "Normal"
$ws.Range("E47").Value = "  +0.00%  "
$ws.Range("E47").Style = "Normal"
$ws.Range("B48").Value = "Arweave"
$ws.Range("B48").Style = "Normal"
$ws.Range("C48").Value = "https://coinranking.com/coin/7XWg41D1+arweave-ar"
$ws.Range("C48").Style = "Normal"
$ws.Range("D48").Value = "'40.15"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -1.26%  "
$ws.Range("E48").Style = "Normal"
$ws.Range("B49").Value = "Monero"
$ws.Range("B49").Style = "Normal"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("C49").Style = "Normal"
$ws.Range("D49").Value = "'141.79"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.47%  "
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "2.788.27"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +1.20%  "
$ws.Range("E50").Style = "Normal"
$ws.Range("D51").Value = "'0.0355"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.63%  "
$ws.Range("E51").Style = "Normal"
